# Change the highlight color of the title text "Ito na ang simula ng
# pagbabago" from yellow to (bright) green, matching the target diff
# where every run in that paragraph moves from
#   <w:highlight w:val="yellow"/>  ->  <w:highlight w:val="green"/>

$d = $word.ActiveDocument

# The title is the first paragraph in the document.
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range

# wdBrightGreen (4) serialises to OOXML w:highlight w:val="green",
# which is exactly the target highlight color in the diff.
$wdBrightGreen = 4
$titleRange.HighlightColorIndex = $wdBrightGreen
